# Auto-generated Excel COM-interop edit script
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4995.25
$ws.Range("I40").Value = 4995.6665
$ws.Range("J40").Value = 4994
$ws.Range("K40").Value = 4995.6665
$ws.Range("L40").Value = 4994
$ws.Range("M40").Value = -4820.6665
$ws.Range("N40").Value = -5344
$ws.Range("H51").Value = 4232.087
$ws.Range("I51").Value = 4498
$ws.Range("J51").Value = 4206.7617
$ws.Range("K51").Value = 4498
$ws.Range("L51").Value = 4206.7617
$ws.Range("M51").Value = -4014
$ws.Range("N51").Value = -5174.7617
$ws.Range("H52").Value = 6439.75
$ws.Range("I52").Value = 5379.5
$ws.Range("K52").Value = 16138.5
$ws.Range("M52").Value = -15978.5
$ws.Range("H69").Value = 10886.667
$ws.Range("I69").Value = 3995
$ws.Range("J69").Value = 14332.5
$ws.Range("K69").Value = 11985
$ws.Range("L69").Value = 42997.5
$ws.Range("M69").Value = -11111
$ws.Range("N69").Value = -44745.5
$ws.Range("H72").Value = 10886.667
$ws.Range("I72").Value = 3995
$ws.Range("J72").Value = 14332.5
$ws.Range("K72").Value = 35955
$ws.Range("L72").Value = 128992.5
$ws.Range("M72").Value = -31587
$ws.Range("N72").Value = -137728.5
$ws.Range("H92").Value = 1341.5
$ws.Range("I92").Value = 1016.6
$ws.Range("K92").Value = 1016.6
$ws.Range("M92").Value = 231.4
$ws.Range("H98").Value = 1960.5714
$ws.Range("I98").Value = 1745.8948
$ws.Range("K98").Value = 1745.8948
$ws.Range("M98").Value = -247.8948
$ws.Range("H107").Value = 1026.0476
$ws.Range("J107").Value = 166
$ws.Range("L107").Value = 166
$ws.Range("N107").Value = -4006
$ws.Range("H122").Value = 1960.5714
$ws.Range("I122").Value = 1745.8948
$ws.Range("K122").Value = 5237.6844
$ws.Range("M122").Value = -2787.6844
$ws.Range("H138").Value = 1957.2354
$ws.Range("J138").Value = 2186
$ws.Range("L138").Value = 6558
$ws.Range("N138").Value = -16838

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1642.3793
$ws.Range("I2").Value = 857.6429000000001
$ws.Range("J2").Value = 2374.8
$ws.Range("K2").Value = 857.6429000000001
$ws.Range("L2").Value = 2374.8
$ws.Range("M2").Value = -744.6429000000001
$ws.Range("N2").Value = -2600.8
$ws.Range("H61").Value = 5299.5
$ws.Range("I61").Value = 2099
$ws.Range("K61").Value = 2099
$ws.Range("M61").Value = -1887
$ws.Range("H63").Value = 1990.2727
$ws.Range("J63").Value = 2357.5715
$ws.Range("L63").Value = 2357.5715
$ws.Range("N63").Value = -3729.5715
$ws.Range("H66").Value = 1990.2727
$ws.Range("J66").Value = 2357.5715
$ws.Range("L66").Value = 11787.8575
$ws.Range("N66").Value = -18651.8575
$ws.Range("H74").Value = 1730.56
$ws.Range("I74").Value = 1667.5
$ws.Range("K74").Value = 1667.5
$ws.Range("M74").Value = -793.5
$ws.Range("H77").Value = 1730.56
$ws.Range("I77").Value = 1667.5
$ws.Range("K77").Value = 8337.5
$ws.Range("M77").Value = -3969.5
$ws.Range("H97").Value = 437.5
$ws.Range("I97").Value = 403.63635
$ws.Range("J97").Value = 561.6667
$ws.Range("K97").Value = 403.63635
$ws.Range("L97").Value = 561.6667
$ws.Range("M97").Value = 92.36365000000001
$ws.Range("N97").Value = -1553.6667
$ws.Range("H116").Value = 1642.3793
$ws.Range("I116").Value = 857.6429000000001
$ws.Range("J116").Value = 2374.8
$ws.Range("K116").Value = 857.6429000000001
$ws.Range("L116").Value = 2374.8
$ws.Range("M116").Value = 1436.3571
$ws.Range("N116").Value = -6962.8
$ws.Range("H122").Value = 3745.2903
$ws.Range("I122").Value = 3778.7307
$ws.Range("K122").Value = 11336.1921
$ws.Range("M122").Value = -8886.1921
$ws.Range("H132").Value = 6953.6587
$ws.Range("I132").Value = 4999.4243
$ws.Range("K132").Value = 14998.2729
$ws.Range("M132").Value = -12468.2729
$ws.Range("H136").Value = 5299.5
$ws.Range("I136").Value = 2099
$ws.Range("K136").Value = 6297
$ws.Range("M136").Value = -3747

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1642.3793
$ws.Range("I3").Value = 857.6429000000001
$ws.Range("J3").Value = 2374.8
$ws.Range("K3").Value = 857.6429000000001
$ws.Range("L3").Value = 2374.8
$ws.Range("M3").Value = -743.6429000000001
$ws.Range("N3").Value = -2602.8
$ws.Range("H64").Value = 1102.4546
$ws.Range("J64").Value = 865.875
$ws.Range("L64").Value = 865.875
$ws.Range("N64").Value = -1315.875
$ws.Range("H67").Value = 1102.4546
$ws.Range("J67").Value = 865.875
$ws.Range("L67").Value = 865.875
$ws.Range("N67").Value = -2425.875
$ws.Range("H86").Value = 2875.05
$ws.Range("I86").Value = 2565.8572
$ws.Range("K86").Value = 2565.8572
$ws.Range("M86").Value = -1442.8572
$ws.Range("H89").Value = 2875.05
$ws.Range("I89").Value = 2565.8572
$ws.Range("K89").Value = 12829.286
$ws.Range("M89").Value = -7213.286
$ws.Range("H134").Value = 4198.5713
$ws.Range("I134").Value = 4198.5713
$ws.Range("K134").Value = 12595.7139
$ws.Range("M134").Value = -10060.7139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1025.125
$ws.Range("I94").Value = 1174
$ws.Range("J94").Value = 876.25
$ws.Range("K94").Value = 1174
$ws.Range("L94").Value = 876.25
$ws.Range("M94").Value = -723
$ws.Range("N94").Value = -1778.25
$ws.Range("H98").Value = 23099.334
$ws.Range("J98").Value = 23099.334
$ws.Range("L98").Value = 23099.334
$ws.Range("N98").Value = -27591.334
$ws.Range("H105").Value = 11228.632
$ws.Range("I105").Value = 16487.916
$ws.Range("J105").Value = 2212.7144
$ws.Range("K105").Value = 16487.916
$ws.Range("L105").Value = 2212.7144
$ws.Range("M105").Value = -14740.916
$ws.Range("N105").Value = -5706.7144
$ws.Range("H134").Value = 3333.3333
$ws.Range("I134").Value = 3300
$ws.Range("K134").Value = 9900
$ws.Range("M134").Value = -7365

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 508.2
$ws.Range("J5").Value = 995.5
$ws.Range("L5").Value = 2986.5
$ws.Range("N5").Value = -3210.5
$ws.Range("H22").Value = 849.5833
$ws.Range("J22").Value = 860.9091
$ws.Range("L22").Value = 2582.7273
$ws.Range("N22").Value = -2920.7273
$ws.Range("H27").Value = 849.5833
$ws.Range("J27").Value = 860.9091
$ws.Range("L27").Value = 2582.7273
$ws.Range("N27").Value = -2786.7273
$ws.Range("H41").Value = 3376.1428
$ws.Range("J41").Value = 3439
$ws.Range("L41").Value = 10317
$ws.Range("N41").Value = -10993
$ws.Range("H50").Value = 687.5
$ws.Range("J50").Value = 415
$ws.Range("L50").Value = 1245
$ws.Range("N50").Value = -2207
$ws.Range("H53").Value = 687.5
$ws.Range("J53").Value = 415
$ws.Range("L53").Value = 1245
$ws.Range("N53").Value = -2207
$ws.Range("H135").Value = 508.2
$ws.Range("J135").Value = 995.5
$ws.Range("L135").Value = 8959.5
$ws.Range("N135").Value = -14029.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3836.0417
$ws.Range("I102").Value = 3730.2727
$ws.Range("K102").Value = 3730.2727
$ws.Range("M102").Value = -2108.2727
$ws.Range("H126").Value = 5737
$ws.Range("I126").Value = 4671.25
$ws.Range("K126").Value = 14013.75
$ws.Range("M126").Value = -11543.75
$ws.Range("H132").Value = 3987.2222
$ws.Range("I132").Value = 4283.5713
$ws.Range("J132").Value = 2950
$ws.Range("K132").Value = 12850.7139
$ws.Range("L132").Value = 8850
$ws.Range("M132").Value = -10320.7139
$ws.Range("N132").Value = -13910

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4928.8184
$ws.Range("I22").Value = 2042.75
$ws.Range("K22").Value = 2042.75
$ws.Range("M22").Value = -1747.75
$ws.Range("H27").Value = 4928.8184
$ws.Range("I27").Value = 2042.75
$ws.Range("K27").Value = 2042.75
$ws.Range("M27").Value = -1935.75
$ws.Range("H40").Value = 7827.3
$ws.Range("I40").Value = 7827.3
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 7827.3
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -7691.3
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3281.8667
$ws.Range("J81").Value = 6439.3335
$ws.Range("L81").Value = 12878.667
$ws.Range("N81").Value = -15000.667
$ws.Range("H84").Value = 3281.8667
$ws.Range("J84").Value = 6439.3335
$ws.Range("L84").Value = 64393.335
$ws.Range("N84").Value = -75001.33499999999
$ws.Range("H113").Value = 950.15
$ws.Range("I113").Value = 927.1429000000001
$ws.Range("J113").Value = 1003.8333
$ws.Range("K113").Value = 2781.4287
$ws.Range("L113").Value = 3011.4999
$ws.Range("M113").Value = -611.4287000000004
$ws.Range("N113").Value = -7351.4999
$ws.Range("H126").Value = 2041.3158
$ws.Range("I126").Value = 2043.6111
$ws.Range("K126").Value = 6130.8333
$ws.Range("M126").Value = -3660.8333
